$wb = $excel.ActiveWorkbook

# ---- Sheet1: Neg_Change ----
$ws1 = $wb.Worksheets.Item("Neg_Change")

$ws1.Cells.Item(2,1).Value = "TCS"
$ws1.Cells.Item(2,2).Value = 3243
$ws1.Cells.Item(2,3).Value = 3259.9
$ws1.Cells.Item(2,4).Value = 3208.6
$ws1.Cells.Item(2,5).Value = 3237.9
$ws1.Cells.Item(2,6).Value = 2544725
$ws1.Cells.Item(2,7).Value = 5244013
$ws1.Cells.Item(2,8).Value = -0.5147370916128545
$ws1.Cells.Item(2,9).Value = "TCS"

$ws1.Cells.Item(3,1).Value = "HINDUNILVR"
$ws1.Cells.Item(3,2).Value = 2350.1
$ws1.Cells.Item(3,3).Value = 2355.1
$ws1.Cells.Item(3,4).Value = 2302
$ws1.Cells.Item(3,5).Value = 2318
$ws1.Cells.Item(3,6).Value = 1898154
$ws1.Cells.Item(3,7).Value = 4008271
$ws1.Cells.Item(3,8).Value = -0.5264407022379475
$ws1.Cells.Item(3,9).Value = "HINDUNILVR"

$ws1.Cells.Item(4,1).Value = "TATASTEEL"
$ws1.Cells.Item(4,2).Value = 167.09
$ws1.Cells.Item(4,3).Value = 168.2
$ws1.Cells.Item(4,4).Value = 163
$ws1.Cells.Item(4,5).Value = 163.5
$ws1.Cells.Item(4,6).Value = 15772400
$ws1.Cells.Item(4,7).Value = 33500436
$ws1.Cells.Item(4,8).Value = -0.5291882171324577
$ws1.Cells.Item(4,9).Value = "TATASTEEL"

$ws1.Cells.Item(5,1).Value = "SUPREMEIND"
$ws1.Cells.Item(5,2).Value = 3352
$ws1.Cells.Item(5,3).Value = 3355.1
$ws1.Cells.Item(5,4).Value = 3310
$ws1.Cells.Item(5,5).Value = 3350
$ws1.Cells.Item(5,6).Value = 279616
$ws1.Cells.Item(5,7).Value = 571062
$ws1.Cells.Item(5,8).Value = -0.5103578945893791
$ws1.Cells.Item(5,9).Value = "SUPREMEIND"

$ws1.Cells.Item(6,1).Value = "GMRAIRPORT"
$ws1.Cells.Item(6,2).Value = 103.51
$ws1.Cells.Item(6,3).Value = 103.74
$ws1.Cells.Item(6,4).Value = 97.93000000000001
$ws1.Cells.Item(6,5).Value = 98.34999999999999
$ws1.Cells.Item(6,6).Value = 15105989
$ws1.Cells.Item(6,7).Value = 30961152
$ws1.Cells.Item(6,8).Value = -0.51209861312654
$ws1.Cells.Item(6,9).Value = "GMRAIRPORT"

$ws1.Cells.Item(7,1).Value = "CAMS"
$ws1.Cells.Item(7,2).Value = 770
$ws1.Cells.Item(7,3).Value = 776.1
$ws1.Cells.Item(7,4).Value = 746.3
$ws1.Cells.Item(7,5).Value = 753.5
$ws1.Cells.Item(7,6).Value = 1513063
$ws1.Cells.Item(7,7).Value = 3467783
$ws1.Cells.Item(7,8).Value = -0.5636800226542433
$ws1.Cells.Item(7,9).Value = "CAMS"

# ---- Sheet2: Pos_Change ----
$ws2 = $wb.Worksheets.Item("Pos_Change")

$ws2.Cells.Item(2,1).Value = "SBILIFE"
$ws2.Cells.Item(2,2).Value = 2005
$ws2.Cells.Item(2,3).Value = 2036
$ws2.Cells.Item(2,4).Value = 2005
$ws2.Cells.Item(2,5).Value = 2023.2
$ws2.Cells.Item(2,6).Value = 1050200
$ws2.Cells.Item(2,7).Value = 749526
$ws2.Cells.Item(2,8).Value = 0.4011521948538143
$ws2.Cells.Item(2,9).Value = "SBILIFE"

$ws2.Cells.Item(3,1).Value = "ICICIBANK"
$ws2.Cells.Item(3,2).Value = 1387
$ws2.Cells.Item(3,3).Value = 1394.9
$ws2.Cells.Item(3,4).Value = 1383.3
$ws2.Cells.Item(3,5).Value = 1387.5
$ws2.Cells.Item(3,6).Value = 7826227
$ws2.Cells.Item(3,7).Value = 4942638
$ws2.Cells.Item(3,8).Value = 0.5834109234785149
$ws2.Cells.Item(3,9).Value = "ICICIBANK"

$ws2.Cells.Item(4,1).Value = "BHARTIARTL"
$ws2.Cells.Item(4,2).Value = 2114.4
$ws2.Cells.Item(4,3).Value = 2114.4
$ws2.Cells.Item(4,4).Value = 2078.9
$ws2.Cells.Item(4,5).Value = 2087
$ws2.Cells.Item(4,6).Value = 4327091
$ws2.Cells.Item(4,7).Value = 2938290
$ws2.Cells.Item(4,8).Value = 0.4726562047993901
$ws2.Cells.Item(4,9).Value = "BHARTIARTL"

$ws2.Cells.Item(5,1).Value = "KOTAKBANK"
$ws2.Cells.Item(5,2).Value = 2158.7
$ws2.Cells.Item(5,3).Value = 2158.9
$ws2.Cells.Item(5,4).Value = 2126
$ws2.Cells.Item(5,5).Value = 2129.8
$ws2.Cells.Item(5,6).Value = 2992047
$ws2.Cells.Item(5,7).Value = 1931695
$ws2.Cells.Item(5,8).Value = 0.548923096037418
$ws2.Cells.Item(5,9).Value = "KOTAKBANK"

$ws2.Cells.Item(6,1).Value = "CIPLA"
$ws2.Cells.Item(6,2).Value = 1515.1
$ws2.Cells.Item(6,3).Value = 1520.5
$ws2.Cells.Item(6,4).Value = 1494.3
$ws2.Cells.Item(6,5).Value = 1497
$ws2.Cells.Item(6,6).Value = 1263559
$ws2.Cells.Item(6,7).Value = 813485
$ws2.Cells.Item(6,8).Value = 0.5532665015335255
$ws2.Cells.Item(6,9).Value = "CIPLA"

$ws2.Cells.Item(7,1).Value = "TRENT"
$ws2.Cells.Item(7,2).Value = 4183.1
$ws2.Cells.Item(7,3).Value = 4202
$ws2.Cells.Item(7,4).Value = 4075.1
$ws2.Cells.Item(7,5).Value = 4080
$ws2.Cells.Item(7,6).Value = 1028704
$ws2.Cells.Item(7,7).Value = 679735
$ws2.Cells.Item(7,8).Value = 0.5133897768983501
$ws2.Cells.Item(7,9).Value = "TRENT"

$ws2.Cells.Item(8,1).Value = "AMBUJACEM"
$ws2.Cells.Item(8,2).Value = 534.7
$ws2.Cells.Item(8,3).Value = 540.55
$ws2.Cells.Item(8,4).Value = 526.55
$ws2.Cells.Item(8,5).Value = 529.9
$ws2.Cells.Item(8,6).Value = 1893815
$ws2.Cells.Item(8,7).Value = 1243514
$ws2.Cells.Item(8,8).Value = 0.5229543052993372
$ws2.Cells.Item(8,9).Value = "AMBUJACEM"

$ws2.Cells.Item(9,1).Value = "IOC"
$ws2.Cells.Item(9,2).Value = 163.66
$ws2.Cells.Item(9,3).Value = 164.42
$ws2.Cells.Item(9,4).Value = 161
$ws2.Cells.Item(9,5).Value = 162.12
$ws2.Cells.Item(9,6).Value = 8744916
$ws2.Cells.Item(9,7).Value = 6124074
$ws2.Cells.Item(9,8).Value = 0.427957271580977
$ws2.Cells.Item(9,9).Value = "IOC"

$ws2.Cells.Item(10,1).Value = "NAUKRI"
$ws2.Cells.Item(10,2).Value = 1386.7
$ws2.Cells.Item(10,3).Value = 1391
$ws2.Cells.Item(10,4).Value = 1361.6
$ws2.Cells.Item(10,5).Value = 1374
$ws2.Cells.Item(10,6).Value = 668131
$ws2.Cells.Item(10,7).Value = 437945
$ws2.Cells.Item(10,8).Value = 0.5256048133898092
$ws2.Cells.Item(10,9).Value = "NAUKRI"

$ws2.Cells.Item(11,1).Value = "TVSMOTOR"
$ws2.Cells.Item(11,2).Value = 3670.1
$ws2.Cells.Item(11,3).Value = 3673.8
$ws2.Cells.Item(11,4).Value = 3590.2
$ws2.Cells.Item(11,5).Value = 3606.1
$ws2.Cells.Item(11,6).Value = 472896
$ws2.Cells.Item(11,7).Value = 331003
$ws2.Cells.Item(11,8).Value = 0.4286758730283411
$ws2.Cells.Item(11,9).Value = "TVSMOTOR"

$ws2.Cells.Item(12,1).Value = "HAVELLS"
$ws2.Cells.Item(12,2).Value = 1432
$ws2.Cells.Item(12,3).Value = 1437
$ws2.Cells.Item(12,4).Value = 1400.1
$ws2.Cells.Item(12,5).Value = 1406
$ws2.Cells.Item(12,6).Value = 642925
$ws2.Cells.Item(12,7).Value = 440018
$ws2.Cells.Item(12,8).Value = 0.4611334081787563
$ws2.Cells.Item(12,9).Value = "HAVELLS"

$ws2.Cells.Item(13,1).Value = "DMART"
$ws2.Cells.Item(13,2).Value = 3950
$ws2.Cells.Item(13,3).Value = 3955.2
$ws2.Cells.Item(13,4).Value = 3826.6
$ws2.Cells.Item(13,5).Value = 3866.9
$ws2.Cells.Item(13,6).Value = 268207
$ws2.Cells.Item(13,7).Value = 184255
$ws2.Cells.Item(13,8).Value = 0.4556294266098613
$ws2.Cells.Item(13,9).Value = "DMART"

$ws2.Cells.Item(14,1).Value = "DLF"
$ws2.Cells.Item(14,2).Value = 719.75
$ws2.Cells.Item(14,3).Value = 720.05
$ws2.Cells.Item(14,4).Value = 685.05
$ws2.Cells.Item(14,5).Value = 686.4
$ws2.Cells.Item(14,6).Value = 3505750
$ws2.Cells.Item(14,7).Value = 2301587
$ws2.Cells.Item(14,8).Value = 0.5231881306246516
$ws2.Cells.Item(14,9).Value = "DLF"

$ws2.Cells.Item(15,1).Value = "SUZLON"
$ws2.Cells.Item(15,2).Value = 52.01
$ws2.Cells.Item(15,3).Value = 53
$ws2.Cells.Item(15,4).Value = 51.57
$ws2.Cells.Item(15,5).Value = 51.69
$ws2.Cells.Item(15,6).Value = 97482417
$ws2.Cells.Item(15,7).Value = 61141969
$ws2.Cells.Item(15,8).Value = 0.5943617550164274
$ws2.Cells.Item(15,9).Value = "SUZLON"

$ws2.Cells.Item(16,1).Value = "NYKAA"
$ws2.Cells.Item(16,2).Value = 255.9
$ws2.Cells.Item(16,3).Value = 259.3
$ws2.Cells.Item(16,4).Value = 249.1
$ws2.Cells.Item(16,5).Value = 253.65
$ws2.Cells.Item(16,6).Value = 5791776
$ws2.Cells.Item(16,7).Value = 3699608
$ws2.Cells.Item(16,8).Value = 0.5655107243794477
$ws2.Cells.Item(16,9).Value = "NYKAA"

$ws2.Cells.Item(17,1).Value = "PHOENIXLTD"
$ws2.Cells.Item(17,2).Value = 1725
$ws2.Cells.Item(17,3).Value = 1770
$ws2.Cells.Item(17,4).Value = 1703.6
$ws2.Cells.Item(17,5).Value = 1713
$ws2.Cells.Item(17,6).Value = 531853
$ws2.Cells.Item(17,7).Value = 354160
$ws2.Cells.Item(17,8).Value = 0.5017308561102327
$ws2.Cells.Item(17,9).Value = "PHOENIXLTD"

$ws2.Cells.Item(18,1).Value = "HINDPETRO"
$ws2.Cells.Item(18,2).Value = 451.8
$ws2.Cells.Item(18,3).Value = 452.95
$ws2.Cells.Item(18,4).Value = 442.7
$ws2.Cells.Item(18,5).Value = 447
$ws2.Cells.Item(18,6).Value = 3238893
$ws2.Cells.Item(18,7).Value = 2114565
$ws2.Cells.Item(18,8).Value = 0.5317065211993957
$ws2.Cells.Item(18,9).Value = "HINDPETRO"

$ws2.Cells.Item(19,1).Value = "HDFCAMC"
$ws2.Cells.Item(19,2).Value = 2578.9
$ws2.Cells.Item(19,3).Value = 2605.1
$ws2.Cells.Item(19,4).Value = 2542.4
$ws2.Cells.Item(19,5).Value = 2555
$ws2.Cells.Item(19,6).Value = 746623
$ws2.Cells.Item(19,7).Value = 471729
$ws2.Cells.Item(19,8).Value = 0.5827371223732269
$ws2.Cells.Item(19,9).Value = "HDFCAMC"

$ws2.Cells.Item(20,1).Value = "OIL"
$ws2.Cells.Item(20,2).Value = 413
$ws2.Cells.Item(20,3).Value = 415
$ws2.Cells.Item(20,4).Value = 401.4
$ws2.Cells.Item(20,5).Value = 404
$ws2.Cells.Item(20,6).Value = 1200785
$ws2.Cells.Item(20,7).Value = 760930
$ws2.Cells.Item(20,8).Value = 0.5780492292326496
$ws2.Cells.Item(20,9).Value = "OIL"

$ws2.Cells.Item(21,1).Value = "IRB"
$ws2.Cells.Item(21,2).Value = 42.7
$ws2.Cells.Item(21,3).Value = 42.79
$ws2.Cells.Item(21,4).Value = 41.51
$ws2.Cells.Item(21,5).Value = 41.68
$ws2.Cells.Item(21,6).Value = 6848620
$ws2.Cells.Item(21,7).Value = 4341458
$ws2.Cells.Item(21,8).Value = 0.577493091030709
$ws2.Cells.Item(21,9).Value = "IRB"

$ws2.Cells.Item(22,1).Value = "BIOCON"
$ws2.Cells.Item(22,2).Value = 395.2
$ws2.Cells.Item(22,3).Value = 395.6
$ws2.Cells.Item(22,4).Value = 381.6
$ws2.Cells.Item(22,5).Value = 381.6
$ws2.Cells.Item(22,6).Value = 6232858
$ws2.Cells.Item(22,7).Value = 4217316
$ws2.Cells.Item(22,8).Value = 0.4779205542103082
$ws2.Cells.Item(22,9).Value = "BIOCON"

$ws2.Cells.Item(23,1).Value = "BANKINDIA"
$ws2.Cells.Item(23,2).Value = 142.35
$ws2.Cells.Item(23,3).Value = 142.89
$ws2.Cells.Item(23,4).Value = 136.77
$ws2.Cells.Item(23,5).Value = 138
$ws2.Cells.Item(23,6).Value = 7093668
$ws2.Cells.Item(23,7).Value = 4937787
$ws2.Cells.Item(23,8).Value = 0.4366087480079639
$ws2.Cells.Item(23,9).Value = "BANKINDIA"

$ws2.Cells.Item(24,1).Value = "SJVN"
$ws2.Cells.Item(24,2).Value = 73.70999999999999
$ws2.Cells.Item(24,3).Value = 74
$ws2.Cells.Item(24,4).Value = 70.7
$ws2.Cells.Item(24,5).Value = 71.2
$ws2.Cells.Item(24,6).Value = 4689166
$ws2.Cells.Item(24,7).Value = 3153252
$ws2.Cells.Item(24,8).Value = 0.4870888847450188
$ws2.Cells.Item(24,9).Value = "SJVN"

$ws2.Cells.Item(25,1).Value = "HUDCO"
$ws2.Cells.Item(25,2).Value = 224
$ws2.Cells.Item(25,3).Value = 224.9
$ws2.Cells.Item(25,4).Value = 210.93
$ws2.Cells.Item(25,5).Value = 212.02
$ws2.Cells.Item(25,6).Value = 6122357
$ws2.Cells.Item(25,7).Value = 3893777
$ws2.Cells.Item(25,8).Value = 0.5723440248375806
$ws2.Cells.Item(25,9).Value = "HUDCO"

$ws2.Cells.Item(26,1).Value = "LAURUSLABS"
$ws2.Cells.Item(26,2).Value = 1026
$ws2.Cells.Item(26,3).Value = 1032
$ws2.Cells.Item(26,4).Value = 996.1
$ws2.Cells.Item(26,5).Value = 1003.8
$ws2.Cells.Item(26,6).Value = 1382140
$ws2.Cells.Item(26,7).Value = 936977
$ws2.Cells.Item(26,8).Value = 0.4751055788989484
$ws2.Cells.Item(26,9).Value = "LAURUSLABS"

$ws2.Cells.Item(27,1).Value = "HFCL"
$ws2.Cells.Item(27,2).Value = 70.98
$ws2.Cells.Item(27,3).Value = 71.33
$ws2.Cells.Item(27,4).Value = 65.7
$ws2.Cells.Item(27,5).Value = 66.3
$ws2.Cells.Item(27,6).Value = 30628278
$ws2.Cells.Item(27,7).Value = 19958374
$ws2.Cells.Item(27,8).Value = 0.5346078793793523
$ws2.Cells.Item(27,9).Value = "HFCL"

$ws2.Cells.Item(28,1).Value = "POONAWALLA"
$ws2.Cells.Item(28,2).Value = 469
$ws2.Cells.Item(28,3).Value = 469.8
$ws2.Cells.Item(28,4).Value = 441.55
$ws2.Cells.Item(28,5).Value = 447.8
$ws2.Cells.Item(28,6).Value = 1484364
$ws2.Cells.Item(28,7).Value = 962629
$ws2.Cells.Item(28,8).Value = 0.5419896969652899
$ws2.Cells.Item(28,9).Value = "POONAWALLA"

